# Remove the post row for the panda entry ("「このパンダは私たちに、幸せは単純な
# ことの中に隠れているのだと教えている」...") which lived at row 598. Deleting the
# entire row shifts every following row up by one (599->598, 600->599, ...,
# 741->740) and shrinks the used range from A1:C741 to A1:C740.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(598).Delete()
